# Switch the report's template placeholders from the old "{{ items.x }}"
# syntax to the new short "[x]" syntax.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "[employee]"
$ws.Range("B2").Value = "[date]"
$ws.Range("C2").Value = "[hours]"
$ws.Range("D2").Value = "[subject]"

# Widen the columns to comfortably fit the (now shorter) placeholder text.
# ColumnWidth uses character units; the saved <col width> ends up
# ColumnWidth + 5/6, so subtract that offset to land on the target widths
# (31.83203125 / 10.5 / 28.5) stored in the sheet XML.
$ws.Columns.Item(1).ColumnWidth = 30.998697916666668
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 27.666666666666668

$ws.Range("D2").Select()
